$d = $word.ActiveDocument

# Locate the paragraph ending in "...change observation date" (Berit's
# comment about wanting today's date to default in the observation date
# field) and position a collapsed range right after it.
$target = "Would prefer to see today’s date as default, and only open calendar if I want to change observation date"
$rng = $d.Content
$rng.Find.Execute($target, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$rng.Collapse(0)
$insertStart = $rng.Start

# Append the reviewer's "– Done" resolution note, matching the pattern used
# elsewhere in the document for items that have been addressed.
$rng.InsertAfter(" – Done")

# wdColor is stored 0x00BBGGRR -> RGB C9211E = 1974729
$redColor = 1974729

# Color each piece separately (touching them as two distinct ranges keeps
# them as two separate runs, " – " and "Done", as in the other instances).
$part1 = $d.Range($insertStart, $insertStart + 3)
$part1.Font.Color = $redColor

$part2 = $d.Range($insertStart + 3, $insertStart + 7)
$part2.Font.Color = $redColor

$word.ActiveDocument.Save()
